# Updated symbol list with latest crypto price/volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of row -> (new Price, new Volume(1h)) values. These are stored as plain
# text in the source data (not numbers), so we force the Text number format on
# each target cell before assigning the value. This preserves exact formatting
# (trailing zeros, no scientific notation, literal "%" suffix) instead of letting
# Excel auto-convert the numeric-looking strings into numbers/percentages.
$updates = @{
    2 = @("288.15", "0.82%")
    3 = @("29.11", "1.27%")
    4 = @("5.278", "4.89%")
    5 = @("0.07013", "4.58%")
    6 = @("7.450", "1.49%")
    7 = @("3.556", "5.18%")
    8 = @("1.396", "2.18%")
    9 = @("0.9034", "-3.92%")
    10 = @("0.1605", "2.59%")
    11 = @("0.07539", "12.20%")
    12 = @("0.07679", "2.15%")
    13 = @("0.02915", "-1.36%")
    14 = @("0.09021", "0.24%")
    15 = @("0.001572", "-1.17%")
    16 = @("0.0006525", "0.63%")
    17 = @("0.006061", "-8.90%")
    18 = @("3.484", "1.13%")
    19 = @("2.232", "-0.67%")
    20 = @("0.3242", "1.10%")
    21 = @("0.1345", "2.75%")
    22 = @("4.005", "-1.86%")
    23 = @("0.1598", "3.19%")
    24 = @("0.04524", "0.88%")
    25 = @("0.001208", "2.55%")
    26 = @("0.004158", "-7.76%")
    27 = @("0.0001167", "-6.48%")
    28 = @("0.0001667", "3.34%")
    40 = @("0.04363", "3.97%")
    41 = @("0.006955", "3.55%")
    42 = @("0.1249", "-0.51%")
    43 = @("0.002064", "2.38%")
    44 = @("0.01160", "-5.15%")
    45 = @("0.00005864", "3.61%")
    47 = @("0.01298", "-0.36%")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dCell = $ws.Cells.Item($row, 4)
    $eCell = $ws.Cells.Item($row, 5)
    $dCell.NumberFormat = "@"
    $eCell.NumberFormat = "@"
    $dCell.Value = $vals[0]
    $eCell.Value = $vals[1]
}
